$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.865.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.256.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.636"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.78"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +14.17%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.02"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.58"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.594.04"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.889"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.259.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.800.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.63"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.89"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.05"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.129"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0776"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.09"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +16.21%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.17"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0323"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.83"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.08"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.56%  "
